# Update "想去人数" (want-to-go count) and "最低票价" (min ticket price)
# figures for the 展览 and 全部类型 sheets, matching the newly scraped data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # F-column (想去人数) updates
    $ws.Range("F3").Value = 7
    $ws.Range("F4").Value = 10546
    $ws.Range("F6").Value = 960
    $ws.Range("F7").Value = 76
    $ws.Range("F8").Value = 1297
    $ws.Range("F9").Value = 7886
    $ws.Range("F11").Value = 456
    $ws.Range("F13").Value = 211
    $ws.Range("F15").Value = 3236
    $ws.Range("F17").Value = 322
    $ws.Range("F18").Value = 718
    $ws.Range("F23").Value = 1685

    # G-column (最低票价) updates: rows 11 and 12 became unavailable for sale
    $ws.Range("G11").Value = "不可售"
    $ws.Range("G12").Value = "不可售"
}
